# Train-loss vs val-loss scaling fix: add cumulative count (W) and
# cumulative percentage (X) columns next to the histogram data in V13:V156.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column W: running cumulative total of column V (W13=V13, Wn = Vn + W(n-1))
$ws.Range("W13").Formula = "=V13"
for ($r = 14; $r -le 156; $r++) {
    $prev = $r - 1
    $ws.Range("W$r").Formula = "=V$r+W$prev"
}
# Match the blue "total" font used elsewhere in the sheet (fontId 3)
$ws.Range("W13:W156").Font.Color = 12611584

# Column X: cumulative percentage of the grand total in V157
for ($r = 13; $r -le 156; $r++) {
    $ws.Range("X$r").Formula = "=W$r/`$V`$157"
}
# Re-use the existing blue Percent style (same one already used by column H)
# instead of inventing a new cellXf.
$ws.Range("H25").Copy() | Out-Null
$ws.Range("X13:X156").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Leave the selection where the author left it after filling the new column.
$ws.Range("X13:X156").Select() | Out-Null
